$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Newly solved LeetCode problem: "Middle of the Linked List" (Easy).
$ws.Range("A9").Value = 876
$ws.Range("B9").Value = "Easy"
$ws.Range("C9").Value = "Middle of the Linked List"
$ws.Range("D9").Value = "http://rb.gy/nrugfa"
$ws.Range("E9").Value = "Pointer"
$ws.Range("F9").Value = "O(n)"
$ws.Range("G9").Value = "Use a slow and a fast pointer."

# Add the external hyperlink on the URL cell, same as the other rows.
$ws.Hyperlinks.Add($ws.Range("D9"), "http://rb.gy/nrugfa")

# Copy the formatting from the row above so the new row matches the
# rest of the table (alignment / wrap / hyperlink style, etc.).
$ws.Range("A8:G8").Copy()
$ws.Range("A9:G9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Reflect the final selection left after editing the sheet.
$ws.Range("E16").Select()
